$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.847.63"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "2.240.54"
$ws.Range("E3").Value = "  -0.10%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.04"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +8.76%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -2.01%  "
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.10"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.16%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.40"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  +18.25%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0975"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -1.98%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.39"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -0.82%  "
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +0.00%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.97"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").Value = "2.572.55"
$ws.Range("E15").Value = "  -0.06%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.11"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +0.49%  "
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.864"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "2.236.18"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").Value = "41.799.44"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "0.0₃0971"
$ws.Range("E20").Value = "  -1.35%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.40"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("E22").Value = "  -0.57%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.27"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +17.58%  "
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.69"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("B25").Value = "WEMIXToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.79"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +3.40%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -0.08%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.52"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +6.71%  "
$ws.Range("E28").Value = "  +2.22%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.73"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.13"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -1.56%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.82"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("E33").Value = "  -0.97%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.51"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +0.73%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0725"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +0.32%  "
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.50"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +19.38%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.70"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("E38").Value = "  +12.97%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0284"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +6.28%  "
$ws.Range("E40").Value = "  +2.27%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.58"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +4.27%  "
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.08"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -0.24%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.214"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +12.33%  "
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.14"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -0.69%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.65"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +15.28%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.97"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +0.85%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.85"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +9.30%  "
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("B50").Value = "BitTorrent-New"
$ws.Range("C50").Value = "https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt"
$ws.Range("D50").Value = "0.0₃0158"
$ws.Range("E50").Value = "  +20.05%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.17"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +8.22%  "
